# Applies the diff: updates the date heading and the 25 division answers
# in the table from the "2023-12-05" worksheet to the "2023-12-06" one.

$d = $word.ActiveDocument

# --- Update the date heading ---
$d.Content.Find.Execute("2023-12-05 Tuesday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2023-12-06 Wednesday", 2)

# --- Update the table of division problems ---
# Only rows 1, 5, 9, 13, 17 of the 20-row table contain data (5 columns each).
# Cell text is set directly (rather than via Find/Replace) because some of
# the original answers are duplicated across cells but map to different new
# values, so positional replacement is required.

$t = $d.Tables.Item(1)

$newValues = @{
    "1,1"  = "90÷9=10, 0"
    "1,2"  = "88÷8=11, 0"
    "1,3"  = "77÷5=15, 2"
    "1,4"  = "82÷2=41, 0"
    "1,5"  = "60÷8=7, 4"

    "5,1"  = "48÷9=5, 3"
    "5,2"  = "40÷8=5, 0"
    "5,3"  = "29÷5=5, 4"
    "5,4"  = "91÷6=15, 1"
    "5,5"  = "99÷9=11, 0"

    "9,1"  = "41÷9=4, 5"
    "9,2"  = "53÷6=8, 5"
    "9,3"  = "29÷9=3, 2"
    "9,4"  = "85÷7=12, 1"
    "9,5"  = "90÷9=10, 0"

    "13,1" = "28÷7=4, 0"
    "13,2" = "67÷8=8, 3"
    "13,3" = "97÷7=13, 6"
    "13,4" = "23÷5=4, 3"
    "13,5" = "17÷4=4, 1"

    "17,1" = "97÷2=48, 1"
    "17,2" = "64÷8=8, 0"
    "17,3" = "69÷5=13, 4"
    "17,4" = "88÷6=14, 4"
    "17,5" = "82÷2=41, 0"
}

foreach ($row in 1, 5, 9, 13, 17) {
    for ($col = 1; $col -le 5; $col++) {
        $key = "$row,$col"
        $t.Cell($row, $col).Range.Text = $newValues[$key]
    }
}
